$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 354; Excel shifts rows 354:430 down to 355:431,
# carrying their existing values/formatting with them.
$ws.Rows("354:354").Insert()

# Populate the newly inserted row 354 with the new record.
$ws.Range("A354").Value = 10
$ws.Range("B354").Value = "Vega Modelo de Temuco"
$ws.Range("C354").Value = "La Araucanía"
$ws.Range("D354").Value = 45211
$ws.Range("D354").NumberFormat = $ws.Range("D355").NumberFormat
$ws.Range("E354").Value = 9
$ws.Range("F354").Value = 100112039
$ws.Range("G354").Value = "Ciboulette"
$ws.Range("H354").Value = "Sin especificar"
$ws.Range("I354").Value = "Primera"
$ws.Range("J354").Value = 80
$ws.Range("K354").Value = 7000
$ws.Range("L354").Value = 7000
$ws.Range("M354").Value = 7000
$ws.Range("N354").Value = "$/docena de atados"
$ws.Range("O354").Value = "Provincia de Cautín"
$ws.Range("P354").Value = 2333
$ws.Range("Q354").Value = 3
$ws.Range("R354").Value = "Hortaliza"
